# Add a new "housing_category" column (big category var) into the
# field-name-validation sheet, inserted as column E (pushing the
# existing "notes" column from E to F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; this shifts the old E ("notes") to F
# and keeps all existing cell contents/formatting intact.
$ws.Columns("E:E").Insert()

# Header
$ws.Range("E1").Value = "housing_category"

# Match the column width Excel ends up with for the new column
# (same width as column D, but not "best fit"/autofit).
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Per-row classification of facility_or_program_type into a broad
# "housing_category" bucket.
$map = @{
    2  = "Permanent Housing"
    3  = "Permanent Housing"
    4  = "Permanent Housing"
    5  = "Permanent Housing"
    6  = "Permanent Housing"
    7  = "Other"
    8  = "Other"
    9  = "Other"
    10 = "Permanent Housing"
    11 = "Permanent Housing"
    12 = "Permanent Housing"
    13 = "Permanent Housing"
    14 = "Other"
    15 = "Permanent Housing"
    16 = "Other"
    17 = "Permanent Housing"
    18 = "Permanent Housing"
    19 = "Permanent Housing"
    20 = "Permanent Housing"
    21 = "Other"
    22 = "Permanent Housing"
    23 = "Permanent Housing"
    24 = "Permanent Housing"
    25 = "Permanent Housing"
    26 = "Permanent Housing"
    27 = "Permanent Housing"
    28 = "Permanent Housing"
    29 = "Permanent Housing"
    30 = "Other"
    31 = "Other"
    32 = "Permanent Housing"
    33 = "Permanent Housing"
    34 = "Permanent Housing"
    35 = "Permanent Housing"
    36 = "Permanent Housing"
    37 = "Permanent Housing"
    38 = "Permanent Housing"
    39 = "Permanent Housing"
    40 = "Other"
    41 = "Other"
    42 = "Other"
    43 = "Other"
    44 = "Other"
    45 = "Other"
    46 = "Other"
    47 = "Other"
    48 = "Other"
    49 = "Permanent Housing"
    50 = "Permanent Housing"
    51 = "Other"
    52 = "Other"
    53 = "Unknown"
    54 = "Other"
}

foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 5).Value = $map[$row]
}

# Restore the selection state to match the final saved workbook
# (the sheet was scrolled down and cell E41 was selected when saved).
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E41").Select()
